$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1704035874439462
$ws.Range("C2").Value = 0.5874439461883408
$ws.Range("J2").Value = 0.04035874439461883
$ws.Range("P2").Value = 0.1076233183856502
$ws.Range("S2").Value = 0.09417040358744394
$ws.Range("B3").Value = 0.007633587786259542
$ws.Range("C3").Value = 0.02290076335877863
$ws.Range("J3").Value = 0.02290076335877863
$ws.Range("P3").Value = 0.7557251908396947
$ws.Range("S3").Value = 0.1908396946564886
$ws.Range("J4").Value = 0.02941176470588235
$ws.Range("P4").Value = 0.7647058823529411
$ws.Range("S4").Value = 0.2058823529411765
$ws.Range("B6").Value = 0.0625
$ws.Range("D6").Value = 0.009615384615384616
$ws.Range("F6").Value = 0.05288461538461538
$ws.Range("J6").Value = 0.2355769230769231
$ws.Range("O6").Value = 0.009615384615384616
$ws.Range("Q6").Value = 0.1682692307692308
$ws.Range("R6").Value = 0.0673076923076923
$ws.Range("S6").Value = 0.3942307692307692
$ws.Range("B7").Value = 0.0963302752293578
$ws.Range("D7").Value = 0.009174311926605505
$ws.Range("F7").Value = 0.06880733944954129
$ws.Range("J7").Value = 0.1422018348623853
$ws.Range("O7").Value = 0.01834862385321101
$ws.Range("Q7").Value = 0.2018348623853211
$ws.Range("R7").Value = 0.07798165137614679
$ws.Range("S7").Value = 0.3853211009174312
$ws.Range("B8").Value = 0.06591337099811675
$ws.Range("D8").Value = 0.01506591337099812
$ws.Range("E8").Value = 0.003766478342749529
$ws.Range("F8").Value = 0.06214689265536723
$ws.Range("J8").Value = 0.1073446327683616
$ws.Range("O8").Value = 0.02259887005649718
$ws.Range("Q8").Value = 0.2259887005649718
$ws.Range("R8").Value = 0.08662900188323917
$ws.Range("S8").Value = 0.4105461393596987
$ws.Range("B9").Value = 0.07096774193548387
$ws.Range("D9").Value = 0.006451612903225806
$ws.Range("F9").Value = 0.1483870967741935
$ws.Range("J9").Value = 0.09677419354838709
$ws.Range("O9").Value = 0.03225806451612903
$ws.Range("Q9").Value = 0.1741935483870968
$ws.Range("R9").Value = 0.09032258064516129
$ws.Range("S9").Value = 0.3806451612903226
$ws.Range("B10").Value = 0.08894878706199461
$ws.Range("D10").Value = 0.0215633423180593
$ws.Range("E10").Value = 0.0008984725965858042
$ws.Range("F10").Value = 0.0637915543575921
$ws.Range("J10").Value = 0.09433962264150944
$ws.Range("O10").Value = 0.01976639712488769
$ws.Range("Q10").Value = 0.252470799640611
$ws.Range("R10").Value = 0.09523809523809523
$ws.Range("S10").Value = 0.3629829290206649
$ws.Range("G11").Value = 0.1584158415841584
$ws.Range("J11").Value = 0.07590759075907591
$ws.Range("K11").Value = 0.1683168316831683
$ws.Range("L11").Value = 0.5907590759075908
$ws.Range("S11").Value = 0.006600660066006601
$ws.Range("G12").Value = 0.772972972972973
$ws.Range("J12").Value = 0.1513513513513514
$ws.Range("K12").Value = 0.01081081081081081
$ws.Range("L12").Value = 0.02702702702702703
$ws.Range("S12").Value = 0.03783783783783784
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2291666666666667
$ws.Range("S13").Value = 0.1041666666666667
$ws.Range("G14").Value = 0.75
$ws.Range("J14").Value = 0.25
$ws.Range("F15").Value = 0.0101010101010101
$ws.Range("H15").Value = 0.2070707070707071
$ws.Range("I15").Value = 0.06060606060606061
$ws.Range("J15").Value = 0.3131313131313131
$ws.Range("K15").Value = 0.0707070707070707
$ws.Range("O15").Value = 0.0303030303030303
$ws.Range("S15").Value = 0.3080808080808081
$ws.Range("F16").Value = 0.01398601398601399
$ws.Range("H16").Value = 0.1818181818181818
$ws.Range("I16").Value = 0.1188811188811189
$ws.Range("J16").Value = 0.3566433566433567
$ws.Range("K16").Value = 0.1188811188811189
$ws.Range("M16").Value = 0.04195804195804196
$ws.Range("O16").Value = 0.04895104895104895
$ws.Range("S16").Value = 0.1188811188811189
$ws.Range("F17").Value = 0.01976284584980237
$ws.Range("H17").Value = 0.2272727272727273
$ws.Range("I17").Value = 0.07114624505928854
$ws.Range("J17").Value = 0.4011857707509882
$ws.Range("K17").Value = 0.1007905138339921
$ws.Range("M17").Value = 0.02371541501976284
$ws.Range("N17").Value = 0.001976284584980237
$ws.Range("O17").Value = 0.0533596837944664
$ws.Range("S17").Value = 0.1007905138339921
$ws.Range("F18").Value = 0.02061855670103093
$ws.Range("H18").Value = 0.2216494845360825
$ws.Range("I18").Value = 0.08762886597938144
$ws.Range("J18").Value = 0.3556701030927835
$ws.Range("K18").Value = 0.1134020618556701
$ws.Range("M18").Value = 0.02061855670103093
$ws.Range("O18").Value = 0.07216494845360824
$ws.Range("S18").Value = 0.1082474226804124
$ws.Range("F19").Value = 0.01265822784810127
$ws.Range("H19").Value = 0.2565400843881857
$ws.Range("I19").Value = 0.06413502109704641
$ws.Range("J19").Value = 0.3493670886075949
$ws.Range("K19").Value = 0.1139240506329114
$ws.Range("M19").Value = 0.02362869198312236
$ws.Range("N19").Value = 0.002531645569620253
$ws.Range("O19").Value = 0.06160337552742616
$ws.Range("S19").Value = 0.1156118143459916
